{"js": "// Applies the three real content edits from the commit:\n//   1. \"... z rozd\u00edlu hodnot je t\u011b\u017ek\u00e9 hodnotit ...\"\n//        -> \"... z rozd\u00edlu hodnot (p\u0159\u00edpadn\u011b z toho, \u017ee jedna je v plusu a\n//            druh\u00e1 taky) je t\u011b\u017ek\u00e9 hodnotit ...\"\n//   2. \"... a z\u00e1vislost tedy t\u011b\u017eko pohledat.\"\n//        -> \"... a z\u00e1vislost tedy t\u011b\u017eko pohledat na prvn\u00ed pohled.\"\n//   3. \"... rozhodn\u011b by mi to usnadnilo orientaci v problematice.\"\n//        -> same text plus a trailing space at the end of the paragraph.\n//\n// (The rest of the diff only wraps already-present words in\n// <w:proofErr w:type=\"gramStart/gramEnd\"/> markers produced by Word's\n// grammar checker on re-save - it does not change any visible text, so\n// there is nothing to reproduce for it here.)\n\nconst body = context.document.body;\n\n// 1) Insert the parenthetical remark about the values' signs.\nconst hit1 = body.search(\"hodnot je t\u011b\u017ek\u00e9 hodnotit\", { matchCase: true });\nawait context.sync();\nif (hit1.items.length === 0) {\n  throw new Error(\"Target text for edit #1 not found\");\n}\nhit1.items[0].insertText(\n  \"hodnot (p\u0159\u00edpadn\u011b z toho, \u017ee jedna je v plusu a druh\u00e1 taky) je t\u011b\u017ek\u00e9 hodnotit\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 2) Append \"na prvn\u00ed pohled\" before the final period.\nconst hit2 = body.search(\"t\u011b\u017eko pohledat.\", { matchCase: true });\nawait context.sync();\nif (hit2.items.length === 0) {\n  throw new Error(\"Target text for edit #2 not found\");\n}\nhit2.items[0].insertText(\"t\u011b\u017eko pohledat na prvn\u00ed pohled.\", \"Replace\");\nawait context.sync();\n\n// 3) Add a trailing space right after \"v problematice.\"\nconst hit3 = body.search(\"v problematice.\", { matchCase: true });\nawait context.sync();\nif (hit3.items.length === 0) {\n  throw new Error(\"Target text for edit #3 not found\");\n}\nhit3.items[0].insertText(\" \", \"After\");\nawait context.sync();\n", "ps1": "# Applies the three real content edits from the commit:\n#   1. \"... z rozdilu hodnot je tezke hodnotit ...\"\n#        -> \"... z rozdilu hodnot (pripadne z toho, ze jedna je v plusu a\n#            druha taky) je tezke hodnotit ...\"\n#   2. \"... a zavislost tedy tezko pohledat.\"\n#        -> \"... a zavislost tedy tezko pohledat na prvni pohled.\"\n#   3. \"... rozhodne by mi to usnadnilo orientaci v problematice.\"\n#        -> same text plus a trailing space at the end of the paragraph.\n#\n# (The rest of the diff only wraps already-present words in\n# <w:proofErr w:type=\"gramStart/gramEnd\"/> markers produced by Word's\n# grammar checker on re-save - it does not change any visible text, so\n# there is nothing to reproduce for it here.)\n\n$d = $word.ActiveDocument\n\n# 1) Insert the parenthetical remark about the values' signs.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"hodnot je t\u011b\u017ek\u00e9 hodnotit\"\n$find.Replacement.Text = \"hodnot (p\u0159\u00edpadn\u011b z toho, \u017ee jedna je v plusu a druh\u00e1 taky) je t\u011b\u017ek\u00e9 hodnotit\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Append \"na prvn\u00ed pohled\" before the final period.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"t\u011b\u017eko pohledat.\"\n$find.Replacement.Text = \"t\u011b\u017eko pohledat na prvn\u00ed pohled.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 3) Add a trailing space right after \"v problematice.\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"v problematice.\"\n$find.Replacement.Text = \"v problematice. \"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
